$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.794.39'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.629.84'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.254'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = '1.854.23'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').Value = '1.623.61'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '0.0₃0759'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.73'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = '25.783.49'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.998'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('E25').Value = '  +1.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.122'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0494'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.903'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Value = '1.140.69'
$ws.Range('E37').Value = '  +2.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.545'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('E39').Value = '  -2.27%  '
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.52'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').Value = '  +0.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.802'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('D46').Value = '1.764.26'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('E48').Value = '  +7.90%  '
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('E51').Value = '  -2.15%  '
